$wb = $excel.ActiveWorkbook

# Add a new "Metadata" worksheet after the last existing sheet
$ws = $wb.Worksheets.Add($null, $wb.Sheets($wb.Sheets.Count))
$ws.Name = "Metadata"

# Populate it with the locale info used by the calc tests
$ws.Range("A1").Value = "Locale"
$ws.Range("B1").Value = "en-GB"
